$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price and Volume columns stay text-formatted (values like "1.00", "0.0978" must not
# be auto-converted to numbers by Excel, matching the original inline-string text content).
$ws.Range("D2:E51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = "59.066.26"
$ws.Range("E2").Value = "  +4.85%  "

# Row 3
$ws.Range("D3").Value = "3.317.26"
$ws.Range("E3").Value = "  +2.68%  "

# Row 4
$ws.Range("E4").Value = "  +0.20%  "

# Row 5
$ws.Range("D5").Value = "407.66"
$ws.Range("E5").Value = "  +2.35%  "

# Row 6
$ws.Range("D6").Value = "110.22"
$ws.Range("E6").Value = "  -1.16%  "

# Row 7
$ws.Range("D7").Value = "0.584"
$ws.Range("E7").Value = "  +4.94%  "

# Row 8
$ws.Range("D8").Value = "1.00"
$ws.Range("E8").Value = "  +0.13%  "

# Row 9
$ws.Range("D9").Value = "0.634"
$ws.Range("E9").Value = "  +2.13%  "

# Row 10
$ws.Range("D10").Value = "39.63"

# Row 11
$ws.Range("D11").Value = "0.0978"
$ws.Range("E11").Value = "  +5.04%  "

# Row 12
$ws.Range("E12").Value = "  +1.25%  "

# Row 13
$ws.Range("D13").Value = "3.847.60"
$ws.Range("E13").Value = "  +2.94%  "

# Row 14
$ws.Range("D14").Value = "8.39"
$ws.Range("E14").Value = "  +3.61%  "

# Row 15
$ws.Range("D15").Value = "19.28"
$ws.Range("E15").Value = "  +0.60%  "

# Row 16
$ws.Range("D16").Value = "3.312.90"
$ws.Range("E16").Value = "  +2.69%  "

# Row 17
$ws.Range("E17").Value = "  -0.48%  "

# Row 18
$ws.Range("D18").Value = "58.984.42"
$ws.Range("E18").Value = "  +4.99%  "

# Row 19
$ws.Range("D19").Value = "10.67"
$ws.Range("E19").Value = "  -3.91%  "

# Row 20
$ws.Range("E20").Value = "  -1.72%  "

# Row 21
$ws.Range("E21").Value = "  +3.88%  "

# Row 22
$ws.Range("D22").Value = "12.90"
$ws.Range("E22").Value = "  -1.35%  "

# Row 23
$ws.Range("D23").Value = "302.73"
$ws.Range("E23").Value = "  +1.47%  "

# Row 24
$ws.Range("D24").Value = "74.80"
$ws.Range("E24").Value = "  -1.62%  "

# Row 25
$ws.Range("E25").Value = "  -0.50%  "

# Row 26
$ws.Range("D26").Value = "28.53"
$ws.Range("E26").Value = "  +1.29%  "

# Row 27
$ws.Range("E27").Value = "  +2.18%  "

# Row 28
$ws.Range("D28").Value = "7.82"
$ws.Range("E28").Value = "  -4.64%  "

# Row 29
$ws.Range("E29").Value = "  -0.74%  "

# Row 30
$ws.Range("D30").Value = "7.26"
$ws.Range("E30").Value = "  -1.90%  "

# Row 31
$ws.Range("D31").Value = "1.00"
$ws.Range("E31").Value = "  -0.01%  "

# Row 32
$ws.Range("E32").Value = "  +1.18%  "

# Row 33
$ws.Range("D33").Value = "11.35"
$ws.Range("E33").Value = "  +1.48%  "

# Row 34
$ws.Range("D34").Value = "40.30"
$ws.Range("E34").Value = "  +8.76%  "

# Row 35
$ws.Range("D35").Value = "0.0524"
$ws.Range("E35").Value = "  +6.78%  "

# Row 36
$ws.Range("E36").Value = "  +0.05%  "

# Row 37
$ws.Range("D37").Value = "51.75"
$ws.Range("E37").Value = "  +0.60%  "

# Row 38
$ws.Range("D38").Value = "3.22"
$ws.Range("E38").Value = "  +3.33%  "

# Row 39
$ws.Range("E39").Value = "  +0.22%  "

# Row 40
$ws.Range("D40").Value = "3.46"
$ws.Range("E40").Value = "  -2.35%  "

# Row 41
$ws.Range("D41").Value = "137.64"
$ws.Range("E41").Value = "  +0.64%  "

# Row 42
$ws.Range("E42").Value = "  +1.81%  "

# Row 43
$ws.Range("D43").Value = "1.89"
$ws.Range("E43").Value = "  -2.13%  "

# Row 45
$ws.Range("E45").Value = "  -5.26%  "

# Row 46
$ws.Range("B46").Value = "TheGraph"
$ws.Range("C46").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D46").Value = "0.277"
$ws.Range("E46").Value = "  -2.19%  "

# Row 47
$ws.Range("B47").Value = "WEMIXToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D47").Value = "2.30"
$ws.Range("E47").Value = "  +9.98%  "

# Row 48
$ws.Range("D48").Value = "22.12"
$ws.Range("E48").Value = "  -1.31%  "

# Row 49
$ws.Range("D49").Value = "2.174.20"
$ws.Range("E49").Value = "  +1.87%  "

# Row 50
$ws.Range("E50").Value = "  +0.07%  "

# Row 51
$ws.Range("B51").Value = "ThetaToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D51").Value = "1.91"
$ws.Range("E51").Value = "  -13.86%  "
